$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.951.28'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.793.87'
$ws.Range("E3").Value = '  -0.85%  '

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '359.11'
$ws.Range("E5").Value = '  +1.19%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.87'
$ws.Range("E6").Value = '  -1.73%  '

# Row 7
$ws.Range("E7").Value = '  -0.68%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.593'
$ws.Range("E9").Value = '  -1.27%  '

# Row 10
$ws.Range("E10").Value = '  -1.41%  '

# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0852'
$ws.Range("E11").Value = '  -0.35%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.134'
$ws.Range("E12").Value = '  +2.24%  '

# Row 13
$ws.Range("E13").Value = '  -1.80%  '

# Row 14
$ws.Range("E14").Value = '  -1.56%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.229.92'
$ws.Range("E15").Value = '  -0.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.792.31'
$ws.Range("E16").Value = '  -0.90%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.945'
$ws.Range("E17").Value = '  +2.74%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.896.80'
$ws.Range("E18").Value = '  +0.28%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.45'
$ws.Range("E19").Value = '  -0.97%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.07'
$ws.Range("E20").Value = '  -2.02%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.39'
$ws.Range("E21").Value = '  +0.25%  '

# Row 22
$ws.Range("E22").Value = '  -1.20%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.27'
$ws.Range("E23").Value = '  +0.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.32'
$ws.Range("E24").Value = '  +0.97%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.77'
$ws.Range("E25").Value = '  -0.47%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.52'
$ws.Range("E26").Value = '  -1.88%  '

# Row 27
$ws.Range("E27").Value = '  +0.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.164'
$ws.Range("E28").Value = '  +17.77%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.28'
$ws.Range("E29").Value = '  -0.18%  '

# Row 30
$ws.Range("E30").Value = '  -3.65%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.18'
$ws.Range("E31").Value = '  +4.89%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '52.16'
$ws.Range("E32").Value = '  -1.11%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.89'
$ws.Range("E33").Value = '  +1.16%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0466'
$ws.Range("E34").Value = '  -2.70%  '

# Row 35
$ws.Range("E35").Value = '  +0.74%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.19'
$ws.Range("E36").Value = '  -3.62%  '

# Row 37
$ws.Range("E37").Value = '  -0.04%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.83'
$ws.Range("E38").Value = '  +2.37%  '

# Row 39
$ws.Range("E39").Value = '  -2.69%  '

# Row 40
$ws.Range("E40").Value = '  -3.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.60'
$ws.Range("E41").Value = '  +1.67%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  -1.63%  '

# Row 43
$ws.Range("E43").Value = '  -1.43%  '

# Row 44
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.14'
$ws.Range("E44").Value = '  -4.04%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.89'
$ws.Range("E45").Value = '  -5.41%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.085.64'
$ws.Range("E46").Value = '  -0.30%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.29'
$ws.Range("E47").Value = '  -1.54%  '

# Row 48
$ws.Range("E48").Value = '  +0.18%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.74'
$ws.Range("E49").Value = '  -3.97%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.949'
$ws.Range("E50").Value = '  -1.75%  '

# Row 51
$ws.Range("E51").Value = '  +30.67%  '
